# CRD falta la habilidad de actualizar los registros
#
# Adds a 4th "solicitud" record block (rows 36-43) to Sheet1, mirroring the
# first record block (rows 7-14): a "Fecha/Para/De" header, the three
# Servicios/Materiales/Equipos checkboxes, a bold column-header row and one
# data row (Cant./Descripcion/Justificacion), merging C:D on the header and
# data rows just like the existing blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Clone the formatting of the first record block onto the new block.
#        Copy+PasteSpecial(formats) row-by-row so each destination row picks
#        up exactly the same cell style the analogous source row uses.
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A36:D36").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:D13").Copy() | Out-Null
$ws.Range("A42:D42").PasteSpecial(-4122) | Out-Null

$ws.Range("A14:D14").Copy() | Out-Null
$ws.Range("A43:D43").PasteSpecial(-4122) | Out-Null

# Row 14's taller custom height doesn't travel with a formats-only paste, so
# match it explicitly on the new data row.
$ws.Rows(43).RowHeight = 27

# --- 2. Fill in the new record's values.
$ws.Range("A36").Value = "Fecha:"
$ws.Range("B36").Value = "02/02/2025"

$ws.Range("A37").Value = "Para:"
$ws.Range("B37").Value = "Prueba"

$ws.Range("A38").Value = "De:"
$ws.Range("B38").Value = "Prueba"

$ws.Range("A39").Value = $false
$ws.Range("B39").Value = "Servicios"

$ws.Range("A40").Value = $false
$ws.Range("B40").Value = "Materiales"

$ws.Range("A41").Value = $true
$ws.Range("B41").Value = "Equipos"

$ws.Range("A42").Value = "Cant."
$ws.Range("B42").Value = "Descripción del material"
$ws.Range("C42").Value = "Justificación"

$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "PruebaPruebaPruebaPruebaPrueba"
$ws.Range("C43").Value = "PruebaPruebaPruebaPrueba"

# --- 3. Merge the Justificacion column across C:D like every other block.
$ws.Range("C42:D42").Merge() | Out-Null
$ws.Range("C43:D43").Merge() | Out-Null
